# Rolling the quarterly "Overview" table forward by one period:
#   - drop the oldest period (column D, "6 ماهه منتهی به 1399/06")
#   - shift every remaining period one column to the left (D<-E, E<-F, ... L<-M)
#   - populate the newly freed last column (M) with the new period's data
#
# Columns D..M are the 10 period columns (column numbers 4..13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$firstCol = 4   # column D
$lastCol  = 13  # column M

function Copy-CellValue($srcRange, $dstRange) {
    # Plain assignment lets Excel's COM layer "smart type" the value (e.g. it
    # will happily reinterpret a bare "1401-10-28" string as a serial date).
    # Re-quoting string values with a leading apostrophe keeps them literal
    # text, matching how the source cell was actually stored.
    $v = $srcRange.Value2
    if ($v -is [string]) {
        $dstRange.Value = "'" + $v
    } else {
        $dstRange.Value = $v
    }
}

function Set-LiteralValue($rng, $v) {
    if ($v -is [string]) {
        $rng.Value = "'" + $v
    } else {
        $rng.Value = $v
    }
}

# Rows that hold one of the 10 rolling period values across D:M, together with
# the brand-new value that belongs in the freshly-opened column M once the
# shift has happened, and any column that isn't a pure left-shift of the old
# data (i.e. a value that was recalculated rather than carried over as-is).
$rows = @(
    @{ Row = 8;  MNew = "12 ماهه منتهی به 1401/12"; Exceptions = @{} },
    @{ Row = 9;  MNew = "1402-02-29";                Exceptions = @{ 9 = "1402-02-29 (8)" } },
    @{ Row = 11; MNew = 9627090;  Exceptions = @{} },
    @{ Row = 12; MNew = -8367351; Exceptions = @{} },
    @{ Row = 13; MNew = 1259740;  Exceptions = @{} },
    @{ Row = 14; MNew = -74158;   Exceptions = @{} },
    @{ Row = 16; MNew = -628;     Exceptions = @{} },
    @{ Row = 17; MNew = 1184953;  Exceptions = @{} },
    @{ Row = 18; MNew = -46506;   Exceptions = @{} },
    @{ Row = 19; MNew = 73347;    Exceptions = @{ 9 = 30552 } },
    @{ Row = 20; MNew = 1211794;  Exceptions = @{ 9 = 738497 } },
    @{ Row = 21; MNew = -129697;  Exceptions = @{} },
    @{ Row = 22; MNew = 1082098;  Exceptions = @{ 9 = 627546 } },
    @{ Row = 24; MNew = 1082098;  Exceptions = @{ 9 = 627546 } },
    @{ Row = 26; MNew = 643259;   Exceptions = @{} }
)

foreach ($item in $rows) {
    $r = $item.Row

    # Shift left: column c picks up whatever was in column c+1 (read it first,
    # since it hasn't been overwritten yet at this point in the loop).
    for ($c = $firstCol; $c -lt $lastCol; $c++) {
        $src = $ws.Cells.Item($r, $c + 1)
        $dst = $ws.Cells.Item($r, $c)
        Copy-CellValue $src $dst
    }

    # The new, newest period lands in the final column (M).
    Set-LiteralValue ($ws.Cells.Item($r, $lastCol)) $item.MNew

    # A few columns were recalculated rather than being a straight carry-over
    # from the next column; patch those in after the shift.
    foreach ($colNum in $item.Exceptions.Keys) {
        Set-LiteralValue ($ws.Cells.Item($r, [int]$colNum)) $item.Exceptions[$colNum]
    }
}
